$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Brief Description" paragraph - merge the two runs that were
# split after "...set up an appointment" / ", or any time in the future..."
# into a single run (Find/Replace naturally coalesces the touched runs).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "or any time in the future, the veteran can choose to", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "or any time in the future, the veteran can choose to", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: "The veteran chooses one of his/her upcoming appointments..."
# paragraph - merge the three runs into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "from the list shown (or, if this is an included use case", $true,
    $false, $false, $false, $false, $true, 1, $false,
    "from the list shown (or, if this is an included use case", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: "INCLUDE [Verify Appointment Status] (...)" paragraph.
#   - merge "INCLUDE [Verify Appointment Status] " + "(" into one run
#   - merge the two REF-field instrText runs into one (" REF VerifyStatusNum \h ")
#     and drop the _GoBack bookmark that used to sit between them
#   - merge ")" + " " into one run
# Field codes (instrText) aren't reachable through Find, so the whole
# paragraph is rewritten in one shot via InsertXML.
# ---------------------------------------------------------------------------
$verifyPara = $d.Paragraphs(27)
$verifyRange = $verifyPara.Range
$verifyXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>INCLUDE [Verify Appointment Status] (</w:t></w:r><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> REF VerifyStatusNum \h </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>2</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r><w:r><w:t xml:space="preserve">) </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$verifyRange.InsertXML($verifyXml) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "Secondary Actors:" paragraph - replace "None" with the new
# description, and re-home the (now orphaned) _GoBack bookmark here, right
# after the new run - mirroring Word's "last edit" bookmark relocation.
# ---------------------------------------------------------------------------
$secondaryPara = $d.Paragraphs(15)
$secondaryRange = $secondaryPara.Range
$secondaryXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>The location/mapping services on the veteran&#8217;s phone</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$secondaryRange.InsertXML($secondaryXml) | Out-Null

Write-Output "done"
